$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Baseline Adj R^2
$ws.Range("B4").Value = -0.03699007383682762
$ws.Range("C4").Value = -0.07258513832132141
$ws.Range("D4").Value = -4.28537385272397
$ws.Range("E4").Value = -35.37974730756679
$ws.Range("F4").Value = -0.1481654852644716
$ws.Range("G4").Value = -0.1596147393180158

# Row 7 - Model R^2
$ws.Range("B7").Value = 0.5207640195710195
$ws.Range("C7").Value = 0.4712270008950245
$ws.Range("D7").Value = 0.4511040903126775
$ws.Range("E7").Value = 0.383217285085094
$ws.Range("F7").Value = 0.6270308418840341
$ws.Range("G7").Value = 0.6381187558940977

# Row 8 - Model Adj R^2
$ws.Range("B8").Value = 0.5082128883787218
$ws.Range("C8").Value = 0.4505376134206387
$ws.Range("D8").Value = -1.744479548436613
$ws.Range("E8").Value = -12.44586318514495
$ws.Range("F8").Value = 0.5724755782496311
$ws.Range("G8").Value = 0.3888760447218207

# Row 9 - Model RMSE
$ws.Range("B9").Value = 2.366492859807684
$ws.Range("C9").Value = 2.405754737303325
$ws.Range("D9").Value = 2.328624938986269
$ws.Range("E9").Value = 2.338093845740941
$ws.Range("F9").Value = 2.312964720183349
$ws.Range("G9").Value = 2.127167318294557

# Row 10 - Model HH
$ws.Range("B10").Value = 64
$ws.Range("C10").Value = 42
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 17
$ws.Range("G10").Value = 1

# Row 11 - Delta R^2
$ws.Range("B11").Value = 0.5312885798669806
$ws.Range("C11").Value = 0.5034251535590623
$ws.Range("D11").Value = 0.5081788608574715
$ws.Range("E11").Value = 1.052013033138616
$ws.Range("F11").Value = 0.6286819584767204
$ws.Range("G11").Value = 0.3247925819216287

# Row 12 - Delta Adj R^2
$ws.Range("B12").Value = 0.5452029622155494
$ws.Range("C12").Value = 0.5231227517419601
$ws.Range("D12").Value = 2.540894304287357
$ws.Range("E12").Value = 22.93388412242184
$ws.Range("F12").Value = 0.7206410635141027
$ws.Range("G12").Value = 0.5484907840398365

# Row 13 - Delta RMSE
$ws.Range("B13").Value = -1.069908851616163
$ws.Range("C13").Value = -0.9554754332054753
$ws.Range("D13").Value = -0.9028991738129926
$ws.Range("E13").Value = -1.507801250450464
$ws.Range("F13").Value = -1.477485418639474
$ws.Range("G13").Value = -0.8030097241382035

# Row 14 - Delta HH
$ws.Range("B14").Value = -10
$ws.Range("C14").Value = -5
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = -6
$ws.Range("G14").Value = 0
